$d = $word.ActiveDocument

# 1) "abre a opção de cadastrar novo tipo de acidente." -> "abre a opção tipo de acidente."
$d.Content.Find.Execute(
    "abre a opção de cadastrar novo ", $true, $false, $false, $false, $false,
    $true, 1, $false, "abre a opção ", 2) | Out-Null

# 2) "O sistema exibe uma nova janela com um formulário a ser preenchido referente ao cadastramento."
#    -> "O sistema exibe uma nova janela, ao clicar no botão novo, um formulário é mostrado referente ao cadastramento."
$d.Content.Find.Execute(
    " janela com um formulário a ser preenchido referente", $true, $false, $false, $false, $false,
    $true, 1, $false, " janela, ao clicar no botão novo, um formulário é mostrado referente", 2) | Out-Null

# 3) "Ao clicar no botão novo, o Usuário ADM. preenche ..." -> "O Usuário ADM. preenche ..."
$d.Content.Find.Execute(
    "Ao clicar no botão novo, o Usuário ADM.", $true, $false, $false, $false, $false,
    $true, 1, $false, "O Usuário ADM.", 2) | Out-Null

# 3b) ", após isso pode clicar no ícone localizado no lado esquerdo do registro para incluir."
#     -> " e clica no ícone de inclusão localizado no lado esquerdo do registro."
$d.Content.Find.Execute(
    ", após isso pode clicar no ícone localizado no lado esquerdo do registro para incluir.", $true, $false, $false, $false, $false,
    $true, 1, $false, " e clica no ícone de inclusão localizado no lado esquerdo do registro.", 2) | Out-Null

# 4) "faz as alterações que desejar e clica em atualizar, a atualização e representada por um ícone localizado no lado esquerdo do registro."
#    -> "faz as alterações que desejar e clica no ícone de atualizar localizado no lado esquerdo do registro."
$d.Content.Find.Execute(
    "e clica em atualizar, a atualização e representada por um ícone localizado", $true, $false, $false, $false, $false,
    $true, 1, $false, "e clica no ícone de atualizar localizado", 2) | Out-Null

# 5) Insert two new empty paragraphs before the (empty) paragraph that holds the
#    "_GoBack" bookmark. That paragraph immediately follows the one ending in
#    "... com sucesso ao Acadsystem." -- locate it via Find (content-addressable,
#    independent of paragraph index / any shifts from the edits above).
$rng = $d.Content
$found = $rng.Find.Execute(
    "Ao final da execução deste caso de uso, o tipo de acidente será adicionado com sucesso ao Acadsystem.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $targetIndex = $i
    }
}

$goBackPara = $d.Paragraphs.Item($targetIndex + 1)
$goBackPara.Range.InsertParagraphBefore()
$goBackPara.Range.InsertParagraphAfter()
